$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet (Spanish -> English) ---
$ws.Name = "Leaf_dmg"

# --- Translate header row (A1:G1) from Spanish to English.
# The engine keeps the table header names and shared-string table in sync
# automatically, and every data cell that referenced the old header strings
# (e.g. the FE/BE/FC/FO treatment codes) is re-pointed to the correct
# (renumbered) shared-string entry as a natural side effect. ---
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Field"
$ws.Range("C1").Value = "Treatment"
$ws.Range("D1").Value = "Repeat"
$ws.Range("E1").Value = "Leaves_dmg_10leaves"
$ws.Range("F1").Value = "Marks_5leaves"
$ws.Range("G1").Value = "Observations"

# --- Narrow columns B:F to fixed (non bestFit/autofit) widths ---
$ws.Columns.Item(2).ColumnWidth = 6.91796875
$ws.Columns.Item(3).ColumnWidth = 11.91796875
$ws.Columns.Item(4).ColumnWidth = 8.91796875
$ws.Columns.Item(5).ColumnWidth = 22.584635416666664
$ws.Columns.Item(6).ColumnWidth = 17.251302083333332

# --- Move/record the active selection on the sheet ---
$ws.Range("G11").Select()
